# Update "想去人数" (interested count) figures on the "展览" and "全部类型"
# worksheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1673
    $ws.Range("F4").Value = 25
    $ws.Range("F6").Value = 453
    $ws.Range("F9").Value = 579
}
